$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 0.01925498268004588
$ws.Cells.Item(2, 4).Value = 0.1667097795782269
$ws.Cells.Item(2, 5).Value = 0.2543256217524998
$ws.Cells.Item(2, 6).Value = 0.8397382084503349
$ws.Cells.Item(2, 7).Value = 0.002417707810360123
$ws.Cells.Item(2, 9).Value = 0.498533902315117
$ws.Cells.Item(2, 10).Value = 0.4359200581151583
$ws.Cells.Item(2, 13).Value = 10.90718762146201
$ws.Cells.Item(2, 15).Value = 2.255514619255877

$ws.Cells.Item(3, 3).Value = 0.01721866860775378
$ws.Cells.Item(3, 4).Value = 0.1633099938318594
$ws.Cells.Item(3, 5).Value = 0.2344801607916764
$ws.Cells.Item(3, 6).Value = 0.8783816250633762
$ws.Cells.Item(3, 7).Value = 0.002422913184348985
$ws.Cells.Item(3, 9).Value = 0.5092947226778008
$ws.Cells.Item(3, 10).Value = 0.392047190437836
$ws.Cells.Item(3, 13).Value = 9.579478518004578
$ws.Cells.Item(3, 15).Value = 2.320578258041337

$ws.Cells.Item(4, 3).Value = 0.0159636772293652
$ws.Cells.Item(4, 4).Value = 0.1613565756093749
$ws.Cells.Item(4, 5).Value = 0.2224950797027887
$ws.Cells.Item(4, 6).Value = 0.9038516372634078
$ws.Cells.Item(4, 7).Value = 0.002426247004729312
$ws.Cells.Item(4, 9).Value = 0.5168163124166583
$ws.Cells.Item(4, 10).Value = 0.3653987628547952
$ws.Cells.Item(4, 13).Value = 8.761117010601197
$ws.Cells.Item(4, 15).Value = 2.364817866070211

$ws.Cells.Item(5, 3).Value = 0.01545111575262581
$ws.Cells.Item(5, 4).Value = 0.1605938705758092
$ws.Cells.Item(5, 5).Value = 0.2176596931779571
$ws.Cells.Item(5, 6).Value = 0.9146618598368406
$ws.Cells.Item(5, 7).Value = 0.002427640371259938
$ws.Cells.Item(5, 9).Value = 0.5201075537421147
$ws.Cells.Item(5, 10).Value = 0.3546085729361721
$ws.Cells.Item(5, 13).Value = 8.426808567558339
$ws.Cells.Item(5, 15).Value = 2.383906797475234

$ws.Cells.Item(6, 3).Value = 0.01536593745467485
$ws.Cells.Item(6, 4).Value = 0.1604692264857448
$ws.Cells.Item(6, 5).Value = 0.2168596713377298
$ws.Cells.Item(6, 6).Value = 0.9164827086960585
$ws.Cells.Item(6, 7).Value = 0.002427873846409188
$ws.Cells.Item(6, 9).Value = 0.5206676091192364
$ws.Cells.Item(6, 10).Value = 0.3528209537351756
$ws.Cells.Item(6, 13).Value = 8.371246646729389
$ws.Cells.Item(6, 15).Value = 2.387140057316799

$ws.Cells.Item(7, 3).Value = 0.01595676922165268
$ws.Cells.Item(7, 4).Value = 0.1613461550328879
$ws.Cells.Item(7, 5).Value = 0.2224296732169151
$ws.Cells.Item(7, 6).Value = 0.9039956921847967
$ws.Cells.Item(7, 7).Value = 0.002426265654903232
$ws.Cells.Item(7, 9).Value = 0.5168597886416535
$ws.Cells.Item(7, 10).Value = 0.3652529666856594
$ws.Cells.Item(7, 13).Value = 8.756611763382296
$ws.Cells.Item(7, 15).Value = 2.365071033906759

$ws.Cells.Item(8, 3).Value = 0.01855385406994259
$ws.Cells.Item(8, 4).Value = 0.165509465036962
$ws.Cells.Item(8, 5).Value = 0.2474403097269544
$ws.Cells.Item(8, 6).Value = 0.8526963598991451
$ws.Cells.Item(8, 7).Value = 0.002419474155459154
$ws.Cells.Item(8, 9).Value = 0.5020521524759261
$ws.Cells.Item(8, 10).Value = 0.4207303133547384
$ws.Cells.Item(8, 13).Value = 10.45002640054662
$ws.Cells.Item(8, 15).Value = 2.277047310750845

$ws.Cells.Item(9, 3).Value = 0.02360840491914473
$ws.Cells.Item(9, 4).Value = 0.1747573082237182
$ws.Cells.Item(9, 5).Value = 0.298159455519297
$ws.Cells.Item(9, 6).Value = 0.7662534958051879
$ws.Cells.Item(9, 7).Value = 0.002407239963772831
$ws.Cells.Item(9, 9).Value = 0.480441802095676
$ws.Cells.Item(9, 10).Value = 0.5320007019089132
$ws.Cells.Item(9, 13).Value = 13.74758621742137
$ws.Cells.Item(9, 15).Value = 2.139282446656637

$ws.Cells.Item(10, 3).Value = 0.02729774720435785
$ws.Cells.Item(10, 4).Value = 0.1822437086135977
$ws.Cells.Item(10, 5).Value = 0.3365754539032508
$ws.Cells.Item(10, 6).Value = 0.7118257097720431
$ws.Cells.Item(10, 7).Value = 0.002398899621352146
$ws.Cells.Item(10, 9).Value = 0.4693320815893358
$ws.Cells.Item(10, 10).Value = 0.6155418225636993
$ws.Cells.Item(10, 13).Value = 16.15910508266336
$ws.Cells.Item(10, 15).Value = 2.060434551546564

$ws.Cells.Item(11, 3).Value = 0.02897078861992952
$ws.Cells.Item(11, 4).Value = 0.1858075070338572
$ws.Cells.Item(11, 5).Value = 0.3543343147341886
$ws.Cells.Item(11, 6).Value = 0.6891369759877435
$ws.Cells.Item(11, 7).Value = 0.002395243255625836
$ws.Cells.Item(11, 9).Value = 0.4653659401093648
$ws.Cells.Item(11, 10).Value = 0.6540031724083804
$ws.Cells.Item(11, 13).Value = 17.25449776207904
$ws.Cells.Item(11, 15).Value = 2.029668998846148

$ws.Cells.Item(12, 3).Value = 0.0296035612972787
$ws.Cells.Item(12, 4).Value = 0.1871804608171317
$ws.Cells.Item(12, 5).Value = 0.3611026623215992
$ws.Cells.Item(12, 6).Value = 0.680851638294854
$ws.Cells.Item(12, 7).Value = 0.002393878253926895
$ws.Cells.Item(12, 9).Value = 0.4640249301223776
$ws.Cells.Item(12, 10).Value = 0.6686392203418166
$ws.Cells.Item(12, 13).Value = 17.66912689178804
$ws.Cells.Item(12, 15).Value = 2.018773600850921

$ws.Cells.Item(13, 3).Value = 0.02946731690455806
$ws.Cells.Item(13, 4).Value = 0.1868837177361229
$ws.Cells.Item(13, 5).Value = 0.3596429995775452
$ws.Cells.Item(13, 6).Value = 0.6826222676281617
$ws.Cells.Item(13, 7).Value = 0.002394171363972672
$ws.Cells.Item(13, 9).Value = 0.4643065131939039
$ws.Cells.Item(13, 10).Value = 0.6654838075427847
$ws.Cells.Item(13, 13).Value = 17.57983572458591
$ws.Cells.Item(13, 15).Value = 2.021086207079492

$ws.Cells.Item(14, 3).Value = 0.02902286275158872
$ws.Cells.Item(14, 4).Value = 0.1859199873445192
$ws.Cells.Item(14, 5).Value = 0.3548902636263875
$ws.Cells.Item(14, 6).Value = 0.6884491455959534
$ws.Cells.Item(14, 7).Value = 0.002395130564742782
$ws.Cells.Item(14, 9).Value = 0.4652523641747024
$ws.Cells.Item(14, 10).Value = 0.6552058204832178
$ws.Cells.Item(14, 13).Value = 17.28861263483253
$ws.Cells.Item(14, 15).Value = 2.028757377249121

$ws.Cells.Item(15, 3).Value = 0.02875052085343555
$ws.Cells.Item(15, 4).Value = 0.185332745576801
$ws.Cells.Item(15, 5).Value = 0.3519848185359109
$ws.Cells.Item(15, 6).Value = 0.6920584373833663
$ws.Cells.Item(15, 7).Value = 0.002395720647068797
$ws.Cells.Item(15, 9).Value = 0.4658528109800528
$ws.Cells.Item(15, 10).Value = 0.6489197529148782
$ws.Cells.Item(15, 13).Value = 17.11020954477289
$ws.Cells.Item(15, 15).Value = 2.03355512748962

$ws.Cells.Item(16, 3).Value = 0.02718830275519224
$ws.Cells.Item(16, 4).Value = 0.1820140455711936
$ws.Cells.Item(16, 5).Value = 0.3354207984262416
$ws.Cells.Item(16, 6).Value = 0.7133508731356315
$ws.Cells.Item(16, 7).Value = 0.002399141325684585
$ws.Cells.Item(16, 9).Value = 0.4696135458955339
$ws.Cells.Item(16, 10).Value = 0.6130379788180846
$ws.Cells.Item(16, 13).Value = 16.0874910489498
$ws.Cells.Item(16, 15).Value = 2.062549624854171

$ws.Cells.Item(17, 3).Value = 0.02622857310684878
$ws.Cells.Item(17, 4).Value = 0.1800191100770689
$ws.Cells.Item(17, 5).Value = 0.3253337289338418
$ws.Cells.Item(17, 6).Value = 0.7269495330124371
$ws.Cells.Item(17, 7).Value = 0.002401274917590041
$ws.Cells.Item(17, 9).Value = 0.472202377440432
$ws.Cells.Item(17, 10).Value = 0.5911469288572846
$ws.Cells.Item(17, 13).Value = 15.45971281459475
$ws.Cells.Item(17, 15).Value = 2.081658604362929

$ws.Cells.Item(18, 3).Value = 0.02567606765265396
$ws.Cells.Item(18, 4).Value = 0.1788865227720038
$ws.Cells.Item(18, 5).Value = 0.3195583977890664
$ws.Cells.Item(18, 6).Value = 0.7349655031113045
$ws.Cells.Item(18, 7).Value = 0.002402515077400532
$ws.Cells.Item(18, 9).Value = 0.4737935141737211
$ws.Cells.Item(18, 10).Value = 0.5785985275455801
$ws.Cells.Item(18, 13).Value = 15.09847496561531
$ws.Cells.Item(18, 15).Value = 2.093128592249116

$ws.Cells.Item(19, 3).Value = 0.02548891444311607
$ws.Cells.Item(19, 4).Value = 0.1785055785216514
$ws.Cells.Item(19, 5).Value = 0.3176074387255881
$ws.Cells.Item(19, 6).Value = 0.7377126983304905
$ws.Cells.Item(19, 7).Value = 0.002402937208829563
$ws.Cells.Item(19, 9).Value = 0.4743496484775989
$ws.Cells.Item(19, 10).Value = 0.574357029801547
$ws.Cells.Item(19, 13).Value = 14.97613742573816
$ws.Cells.Item(19, 15).Value = 2.097093768175995

$ws.Cells.Item(20, 3).Value = 0.02633078925622101
$ws.Cells.Item(20, 4).Value = 0.1802299330772712
$ws.Cells.Item(20, 5).Value = 0.3264047545366822
$ws.Cells.Item(20, 6).Value = 0.7254817529642921
$ws.Cells.Item(20, 7).Value = 0.002401046451881321
$ws.Cells.Item(20, 9).Value = 0.4719161919145449
$ws.Cells.Item(20, 10).Value = 0.5934728038227206
$ws.Cells.Item(20, 13).Value = 15.52655664789017
$ws.Cells.Item(20, 15).Value = 2.079574694085551

$ws.Cells.Item(21, 3).Value = 0.0291534307581145
$ws.Cells.Item(21, 4).Value = 0.1862024170652887
$ws.Cells.Item(21, 5).Value = 0.3562850554305186
$ws.Cells.Item(21, 6).Value = 0.6867292651330246
$ws.Cells.Item(21, 7).Value = 0.002394848294481388
$ws.Cells.Item(21, 9).Value = 0.4649701424447272
$ws.Cells.Item(21, 10).Value = 0.6582227246862544
$ws.Cells.Item(21, 13).Value = 17.37415607949436
$ws.Cells.Item(21, 15).Value = 2.026483511939318

$ws.Cells.Item(22, 3).Value = 0.0309936875341208
$ws.Cells.Item(22, 4).Value = 0.1902426449148038
$ws.Cells.Item(22, 5).Value = 0.3760682357732605
$ws.Cells.Item(22, 6).Value = 0.66319231720842
$ws.Cells.Item(22, 7).Value = 0.002390911488752711
$ws.Cells.Item(22, 9).Value = 0.4613703824908484
$ws.Cells.Item(22, 10).Value = 0.7009606292401713
$ws.Cells.Item(22, 13).Value = 18.58069721746864
$ws.Cells.Item(22, 15).Value = 1.996195310636125

$ws.Cells.Item(23, 3).Value = 0.03001192381132967
$ws.Cells.Item(23, 4).Value = 0.1880735399950595
$ws.Cells.Item(23, 5).Value = 0.365485326812319
$ws.Cells.Item(23, 6).Value = 0.6755877380787396
$ws.Cells.Item(23, 7).Value = 0.002393002273892644
$ws.Cells.Item(23, 9).Value = 0.4632041411190997
$ws.Cells.Item(23, 10).Value = 0.6781101816841897
$ws.Cells.Item(23, 13).Value = 17.93681085018397
$ws.Cells.Item(23, 15).Value = 2.011950096546173

$ws.Cells.Item(24, 3).Value = 0.02628457965195707
$ws.Cells.Item(24, 4).Value = 0.1801345754219028
$ws.Cells.Item(24, 5).Value = 0.3259204697484392
$ws.Cells.Item(24, 6).Value = 0.7261447204997111
$ws.Cells.Item(24, 7).Value = 0.002401149699130844
$ws.Cells.Item(24, 9).Value = 0.4720452565817865
$ws.Cells.Item(24, 10).Value = 0.5924211607592724
$ws.Cells.Item(24, 13).Value = 15.49633755140241
$ws.Cells.Item(24, 15).Value = 2.080515324019444

$ws.Cells.Item(25, 3).Value = 0.02224522436744536
$ws.Cells.Item(25, 4).Value = 0.1721364408993225
$ws.Cells.Item(25, 5).Value = 0.2842463651525264
$ws.Cells.Item(25, 6).Value = 0.7880769571202038
$ws.Cells.Item(25, 7).Value = 0.002410434881906519
$ws.Cells.Item(25, 9).Value = 0.4854683230602106
$ws.Cells.Item(25, 10).Value = 0.5016060845674133
$ws.Cells.Item(25, 13).Value = 12.85770966492498
$ws.Cells.Item(25, 15).Value = 2.17270578697125
